$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10_27")

# New header cell + z-score column (G) for GC measurements block
$ws.Range("G11").Value = "z-score"
$ws.Range("G12").Formula = '=(F12-$F$18)/$F$19'
$ws.Range("G13:G16").Formula = '=(F13-$F$18)/$F$19'

# New summary rows for flow rate / flux
$ws.Range("A20").Value = "Flow rate (mmol/hr)"
$ws.Range("B20").Formula = '=B17*60/1000/22.4'
$ws.Range("A21").Value = "Flux (mmol/g(DCW)/hr"
$ws.Range("B21").Formula = '=B20/B1'

# Restore cursor/selection position on each sheet, matching where the
# author last left the cursor before saving
$ws1 = $wb.Worksheets.Item("10_22")
[void]$ws1.Range("F18").Select()

[void]$ws.Range("F28").Select()
$ws.Activate()
